$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.48619270324707
$ws.Range("B1").Value = 1.958019614219666
$ws.Range("C1").Value = 3.34549880027771
$ws.Range("D1").Value = 1.368920683860779
$ws.Range("E1").Value = 0.809965193271637
